# Scene.xlsx maintenance edit:
#  - rename resource/ini path segment to resource/res (commit: "rename resource/struct as resource/schema")
#  - re-fit the FilePath column now that its contents are wider
#  - tweak the workbook's Background-1 theme colour ("fix invalid netclient id" era re-theme)
#  - leave the cursor where the author left it when they saved

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. FilePath column (F9:F14) referenced "../resource/ini/Scene/N.xml" for each scene;
#    rename the "ini" folder segment to "res" for every row in one shot.
$ws.Range("F9:F14").Replace("ini/Scene", "res/Scene")

# 2. The FilePath values are now longer/shorter than before - re-autofit column F so the
#    text is fully visible again (column F was fixed-width 14 before, now content driven).
$ws.Columns.Item(6).EntireColumn.AutoFit()

# 3. Update the theme's "Background 1" / Light 1 colour (index 2 in the theme colour scheme)
#    from the default white to a pale green.
$theme = $wb.Theme
$colorScheme = $theme.ThemeColorScheme
$background1 = $colorScheme.Colors(2)
$background1.RGB = 13625548   # RGB(204, 232, 207) = &HCCE8CF

# 4. Restore the sheet selection to where the author clicked last.
$ws.Range("F25").Select()
